$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title (paragraph 1)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Echoes of the Past: Historical Lessons for Modern Governance",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "A Glimpse into the Art of Government: An Exploration of Civics and Politics", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Author name (paragraph 2)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Marcus Walton",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "Clara Bennett", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. E-mail address (paragraph 3): marcuswalton56@abromail.net -> clara.bennett88@institute.edu
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3start = $p3.Range.Start
$p3textEnd = $p3.Range.End - 1
$emailRng = $d.Range($p3start, $p3textEnd)
$emailRng.Text = "clara.bennett88@institute.edu"
$emailRng2 = $d.Range($p3start, $p3start + 30)
$emailRng2.Font.Name = "Times New Roman"
$emailRng2.Font.NameAscii = "Times New Roman"
$emailRng2.Font.Size = 16
$emailRng2.Font.Color = 0

# ---------------------------------------------------------------------
# 4. Main body paragraph (paragraph 5): trim to first three sentences, then
#    replace their wording with the new civics/politics text.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$fullText = $p5.Range.Text
$cutIdx = $fullText.IndexOf("leadership.") + "leadership.".Length
$p5start = $p5.Range.Start
$keepEnd = $p5start + $cutIdx
$p5paraEnd = $p5.Range.End - 1
$tailRng = $d.Range($keepEnd, $p5paraEnd)
$tailRng.Text = ""

$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute(
    "History, like a murmuring river, whispers tales of triumphs and tribulations, inviting us to glean wisdom from the annals of time",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "The exploration of government and politics provides a lens through which we can examine the interplay of power dynamics, decision-making processes, and the quest for justice", 2) | Out-Null

$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute(
    " It is a tapestry woven with threads of human experience, where patterns emerge, offering guidance for the complexities of modern governance",
    $true, $true, $false, $false, $false, $true, 1, $false,
    " It encourages us to think critically, to challenge assumptions, and to recognize the interconnections between our actions and their broader implications", 2) | Out-Null

$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute(
    " By delving into the chronicles of bygone eras, we can unearth valuable lessons that illuminate paths towards effective and just leadership",
    $true, $true, $false, $false, $false, $true, 1, $false,
    " As we navigate the complexities of governance and political engagement, we gain a deeper appreciation for the rights and responsibilities that come with being a citizen, and we embrace the opportunity to contribute to a better future for ourselves and for generations to come", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Summary paragraph (paragraph 7)
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Find.Execute(
    "History, as a mirror to the present, holds a wealth of lessons for modern governance",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "In conclusion, government and politics are fundamental pillars of human society, shaping the structures, processes, and relationships that define how we live together", 2) | Out-Null

$p7 = $d.Paragraphs.Item(7)
$p7.Range.Find.Execute(
    " By studying the echoes of the past - the triumphs and tribulations of civilizations, the struggles for justice and equality, and the intricacies of diplomacy and statecraft - leaders can gain insights into the challenges they face and the paths they must tread",
    $true, $true, $false, $false, $false, $true, 1, $false,
    " The study of government and politics provides a critical lens through which we can examine the intricacies of governance, the interplay of power, and the quest for a just and equitable society", 2) | Out-Null

$p7 = $d.Paragraphs.Item(7)
$p7.Range.Find.Execute(
    " History serves as a constant reminder of the importance of balance, empathy, and collective action in shaping a better world, where the lessons of the past illuminate the road ahead",
    $true, $true, $false, $false, $false, $true, 1, $false,
    " It equips us with the knowledge, skills, and values necessary to navigate the political landscape, to participate effectively in the decision-making process, and to work towards a better future for all", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Append a new, fully empty paragraph right before the section break.
# ---------------------------------------------------------------------
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")

# ---------------------------------------------------------------------
# 7. Fix the misspelled font name everywhere in the body ("TimesNewToman"
#    -> "Times New Roman"). Apply per-paragraph (excluding the paragraph
#    mark) so no stray formatting is added to the paragraph marks.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End
    if ($pEnd -gt $pStart) {
        $pEnd = $pEnd - 1
    }
    if ($pEnd -gt $pStart) {
        $fr = $d.Range($pStart, $pEnd)
        $fr.Font.Name = "Times New Roman"
        $fr.Font.NameAscii = "Times New Roman"
    }
}
